$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SIML")

# Income Statement section
$ws.Range("D14").Value = 1100
$ws.Range("D17").Value = 2100
$ws.Range("D24").Value = 0
$ws.Range("D26").Value = -4400
$ws.Range("D27").Value = -4400
$ws.Range("D33").Value = -4400
$ws.Range("D35").Value = -4400

# Balance Sheet section
$ws.Range("D58").Value = 800
$ws.Range("D81").Value = -4400

# Cash Flow Statement section
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("F91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("D94").Value = 0
